$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# New data rows to append (Fecha serial, Pruebas Realizadas, Pruebas Positivas,
# Clinicamente Estables, Clinicamente Graves, Cuidados Intensivos)
$newData = @(
    @(44092, 1204, 491, 607, 98, 25),
    @(44093, 1295, 532, 641, 104, 22),
    @(44094, 1163, 473, 661, 111, 21)
)

foreach ($rowVals in $newData) {
    $lastRow = $lo.Range.Rows.Count + $lo.Range.Row - 1

    # Grow the table by one row so the table/autofilter ref and dimension extend.
    $lo.ListRows.Add() | Out-Null

    $newRowIndex = $lastRow + 1

    # Copy formatting (number formats/styles) from the row above into the new row.
    $ws.Range("A$lastRow`:F$lastRow").Copy()
    $ws.Range("A$newRowIndex").PasteSpecial(-4122)

    $ws.Cells.Item($newRowIndex, 1).Value = $rowVals[0]
    $ws.Cells.Item($newRowIndex, 2).Value = $rowVals[1]
    $ws.Cells.Item($newRowIndex, 3).Value = $rowVals[2]
    $ws.Cells.Item($newRowIndex, 4).Value = $rowVals[3]
    $ws.Cells.Item($newRowIndex, 5).Value = $rowVals[4]
    $ws.Cells.Item($newRowIndex, 6).Value = $rowVals[5]
}

$ws.Application.CutCopyMode = $false

# Match the final selection left by the editor after entering the last value.
$ws.Range("F192").Select()
